# "Adicao dados da turma" - add class/group info to the title slide's
# DISCIPLINA textbox, and refresh the cached date / slide-number
# placeholder fields that PowerPoint re-stamped on the layouts + notes
# master when the deck was re-saved.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1: append a new "3SI." line under "COMPLIANCE & QUALITY
#    ASSURANCE" in the DISCIPLINA textbox (shape keeps its auto-fit,
#    so the box grows to the new height on its own).
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.HasTextFrame) {
        $full = $shp.TextFrame.TextRange.Text
        if ($full -like "*COMPLIANCE*QUALITY ASSURANCE*") {
            $shp.TextFrame.TextRange.InsertAfter([char]13 + "3SI.")
        }
    }
}

# ---------------------------------------------------------------------
# 2) Refresh the cached "datetimeFigureOut" and "slidenum" placeholder
#    text on every slide layout and on the notes master.
# ---------------------------------------------------------------------
$oldDate = "10/02/2020"
$newDate = "05/03/2020"
$oldNum  = [char]0x2039 + "n" + [char]0x00BA + [char]0x203A
$newNum  = [char]0x2039 + "#" + [char]0x203A

function Update-FieldText($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $t = $sh.TextFrame.TextRange.Text
            if ($t -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            } elseif ($t -eq $oldNum) {
                $sh.TextFrame.TextRange.Text = $newNum
            }
        }
    }
}

$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    Update-FieldText $layouts.Item($li).Shapes
}

Update-FieldText $p.NotesMaster.Shapes
